$d = $word.ActiveDocument

# --- 1) Remove the "Meta description: ..." paragraph that currently sits
#        right under the title heading. ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- 2) At the end of the document, the paragraph that holds the AI
#        image-generation prompt ("Please create a cartoon-style feature
#        image...") is turned into two paragraphs:
#          a) a new bold paragraph repeating the page title
#          b) the same paragraph (still italic), but its text is replaced
#             with the meta description that used to live at the top.
$total = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($total)

# Target just the paragraph's own text (exclude the trailing paragraph
# mark) so the inserted XML below produces exactly two paragraphs instead
# of leaving a stray empty one behind.
$targetRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Florageddon! for Free - A Thrilling Slot Game</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Find out why Florageddon! is a must-play online slot game. Read our review and play for free to experience the thrill!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xml)
